$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Germany) - D2,E2,F2 set to 0; G2,H2 cleared entirely
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 0
$ws.Range("G2:H2").ClearContents()

# Row 3 (USA_NYC)
$ws.Range("D3").Value = -0.005516924598505269
$ws.Range("E3").Value = 0.00005403861845915836
$ws.Range("F3").Value = -0.005570963216964427
$ws.Range("G3").Value = 0.009606862369154949
$ws.Range("H3").Value = 0.9903931376308451

# Row 4 (SouthKorea)
$ws.Range("D4").Value = -0.009464650821448846
$ws.Range("E4").Value = 0.0004602093717535852
$ws.Range("F4").Value = -0.009924860193202432
$ws.Range("G4").Value = 0.0443145198859854
$ws.Range("H4").Value = 0.9556854801140146

# Row 5 (US)
$ws.Range("D5").Value = -0.01056873432908085
$ws.Range("E5").Value = -0.00640649477880307
$ws.Range("F5").Value = -0.004162239550277782
$ws.Range("G5").Value = 0.6061742664090822
$ws.Range("H5").Value = 0.3938257335909177

# Row 6 (China)
$ws.Range("D6").Value = -0.01601114141880398
$ws.Range("E6").Value = -0.0005175320192456988
$ws.Range("F6").Value = -0.01549360939955828
$ws.Range("G6").Value = 0.03232324327845191
$ws.Range("H6").Value = 0.9676767567215481

# Row 7 (France)
$ws.Range("D7").Value = -0.03308171270326781
$ws.Range("E7").Value = -0.01507069312187183
$ws.Range("F7").Value = -0.01801101958139599
$ws.Range("G7").Value = 0.4555596397638488
$ws.Range("H7").Value = 0.5444403602361511

# Row 8 (USA_WA)
$ws.Range("D8").Value = -0.03674145150616512
$ws.Range("E8").Value = -0.0142996188151363
$ws.Range("F8").Value = -0.02244183269102882
$ws.Range("G8").Value = 0.3891958055259972
$ws.Range("H8").Value = 0.6108041944740028

# Row 9 (Spain)
$ws.Range("D9").Value = -0.07940112835449982
$ws.Range("E9").Value = -0.03089001912198475
$ws.Range("F9").Value = -0.04851110923251507
$ws.Range("G9").Value = 0.3890375333719569
$ws.Range("H9").Value = 0.6109624666280432

# Row 10 (Italy)
$ws.Range("D10").Value = -0.0995778431740493
$ws.Range("E10").Value = -0.04091651535933921
$ws.Range("F10").Value = -0.05866132781471008
$ws.Range("G10").Value = 0.4108997951263355
$ws.Range("H10").Value = 0.5891002048736646
